# Update the "2.Distance Measurement" data sheet with the new high
# resolution test data captured on laptop (61 rows x 2 cols: A = distance, B = value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.84, 0.3834092881944444),
    @(1.17, 0.5770399305555556),
    @(1.35, 0.5110134548611112),
    @(2.47, 0.2734375),
    @(0.6899999999999999, 0.9430338541666666),
    @(0.65, 0.9764539930555556),
    @(2.28, 0.3023003472222222),
    @(4.09, 0.09136284722222222),
    @(1.08, 0.6234266493055556),
    @(2.81, 0.2265082465277778),
    @(2.88, 0.2183702256944444),
    @(2.87, 0.2192925347222222),
    @(2.83, 0.2242838541666667),
    @(2.84, 0.2235243055555556),
    @(2.83, 0.2245551215277778),
    @(2.9, 0.2160373263888889),
    @(2.87, 0.2194010416666667),
    @(2.88, 0.21875),
    @(2.86, 0.22119140625),
    @(2.84, 0.2231987847222222),
    @(2.8, 0.2284071180555556),
    @(2.86, 0.2206488715277778),
    @(2.87, 0.21923828125),
    @(2.85, 0.2221137152777778),
    @(5.33, -0.0078125),
    @(5.33, -0.0078125),
    @(1.54, 0.4543728298611111),
    @(5.22, 0),
    @(2.39, 0.2853732638888889),
    @(2.19, 0.3167860243055556),
    @(0.63, 0.9921875),
    @(3.23, 0.1764322916666667),
    @(0.63, 0.9921875),
    @(0.63, 0.9921875),
    @(0.63, 0.9921875),
    @(1.27, 0.53759765625),
    @(1.7, 0.4136827256944444),
    @(2.11, 0.3317057291666667),
    @(1.13, 0.5975477430555556),
    @(3.71, 0.1270073784722222),
    @(1.2, 0.5674370659722222),
    @(1, 0.6704644097222222),
    @(1, 0.6674262152777778),
    @(1.01, 0.662109375),
    @(0.99, 0.6754557291666666),
    @(0.99, 0.6780056423611112),
    @(1, 0.6688910590277778),
    @(1.01, 0.6654730902777778),
    @(5.22, 0),
    @(0.93, 0.7210286458333334),
    @(0.92, 0.7255316840277778),
    @(0.93, 0.7190212673611112),
    @(0.78, 0.8559027777777778),
    @(0.79, 0.8513454861111112),
    @(0.84, 0.7976345486111112),
    @(3.17, 0.1841362847222222),
    @(0.8100000000000001, 0.8272026909722222),
    @(0.79, 0.8484700520833334),
    @(0.8, 0.8336046006944444),
    @(0.79, 0.8463541666666666),
    @(0.8100000000000001, 0.8326822916666666)
)

$numRows = 61
$numCols = 2

$arr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $arr[$i, 0] = $data[$i][0]
    $arr[$i, 1] = $data[$i][1]
}

$rng = $ws.Range("A1:B$numRows")
$rng.Value = $arr
